$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.091.93"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "2.431.71"
$ws.Range("E3").Value = "  +4.27%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.12"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.60"
$ws.Range("E6").Value = "  +5.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").Value = "2.429.20"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("E11").Value = "  +4.52%  "
$ws.Range("E13").Value = "  +3.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.75"
$ws.Range("E14").Value = "  +8.88%  "
$ws.Range("D15").Value = "2.855.62"
$ws.Range("E15").Value = "  +3.89%  "
$ws.Range("D16").Value = "62.006.25"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("E17").Value = "  +5.19%  "
$ws.Range("D18").Value = "2.425.84"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.31"
$ws.Range("E20").Value = "  +9.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.23"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.83"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.03"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.52"
$ws.Range("E27").Value = "  +10.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("E29").Value = "  +13.38%  "
$ws.Range("D30").Value = "0.0₃0789"
$ws.Range("E30").Value = "  +7.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  +4.47%  "
$ws.Range("E32").Value = "  +6.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "171.10"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("E34").Value = "  +4.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.396"
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "374.39"
$ws.Range("E36").Value = "  +16.15%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.56"
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").Value = "  +9.87%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +9.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.10"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.78"
$ws.Range("E43").Value = "  +6.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.67"
$ws.Range("E44").Value = "  +5.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.73"
$ws.Range("E45").Value = "  +8.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0957"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.589"
$ws.Range("E47").Value = "  +4.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  +5.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.98"
$ws.Range("E49").Value = "  +6.37%  "
$ws.Range("E50").Value = "  +3.66%  "
$ws.Range("E51").Value = "  +11.74%  "
